$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Assign a value while forcing text storage so Excel does not
    # auto-convert numeric-looking strings (e.g. "505.86") into
    # real numbers, and restore the cell style so no new style is
    # introduced for the cell.
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

$ws.Range('D2').Value = '59.354.62'
$ws.Range('E2').Value = '  -3.60%  '
$ws.Range('D3').Value = '2.843.94'
$ws.Range('E3').Value = '  -4.30%  '
$ws.Range('E4').Value = '  -0.10%  '
Set-TextValue 'D5' '505.86'
$ws.Range('E5').Value = '  -6.32%  '
Set-TextValue 'D6' '136.87'
$ws.Range('E6').Value = '  -8.49%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  -6.23%  '
$ws.Range('D9').Value = '2.844.55'
$ws.Range('E9').Value = '  -4.69%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D10' '0.104'
$ws.Range('E10').Value = '  -7.85%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D11' '5.98'
$ws.Range('E11').Value = '  -2.27%  '
Set-TextValue 'D12' '0.348'
$ws.Range('E12').Value = '  -4.72%  '
$ws.Range('D13').Value = '3.346.36'
$ws.Range('E13').Value = '  -4.32%  '
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('D15').Value = '59.420.44'
$ws.Range('E15').Value = '  -3.63%  '
Set-TextValue 'D16' '21.80'
$ws.Range('E16').Value = '  -8.23%  '
$ws.Range('D17').Value = '2.843.70'
$ws.Range('E17').Value = '  -4.97%  '
Set-TextValue 'D18' '0.0000136'
$ws.Range('E18').Value = '  -7.01%  '
$ws.Range('E19').Value = '  -7.07%  '
Set-TextValue 'D20' '11.13'
$ws.Range('E20').Value = '  -6.82%  '
Set-TextValue 'D21' '352.78'
$ws.Range('E21').Value = '  -5.70%  '
Set-TextValue 'D22' '6.30'
$ws.Range('E22').Value = '  -5.62%  '
Set-TextValue 'D23' '0.998'
$ws.Range('E23').Value = '  -0.14%  '
Set-TextValue 'D24' '5.64'
$ws.Range('E24').Value = '  -0.42%  '
Set-TextValue 'D25' '63.18'
$ws.Range('E25').Value = '  -3.85%  '
$ws.Range('E26').Value = '  -7.83%  '
$ws.Range('E27').Value = '  -8.30%  '
Set-TextValue 'D28' '0.999'
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('E29').Value = '  -7.73%  '
$ws.Range('D30').Value = '0.0₃0822'
$ws.Range('E30').Value = '  -9.39%  '
Set-TextValue 'D31' '1.00'
$ws.Range('E31').Value = '  +0.04%  '
Set-TextValue 'D32' '1.62'
$ws.Range('E32').Value = '  -5.64%  '
Set-TextValue 'D33' '19.09'
$ws.Range('E33').Value = '  -5.99%  '
Set-TextValue 'D34' '150.70'
$ws.Range('E34').Value = '  -5.59%  '
$ws.Range('E35').Value = '  -7.41%  '
$ws.Range('E36').Value = '  -7.80%  '
Set-TextValue 'D37' '0.938'
$ws.Range('E37').Value = '  -11.01%  '
$ws.Range('E38').Value = '  -8.62%  '
Set-TextValue 'D39' '36.55'
$ws.Range('E39').Value = '  -1.80%  '
$ws.Range('D40').Value = '2.231.88'
$ws.Range('E40').Value = '  -7.23%  '
Set-TextValue 'D41' '0.632'
$ws.Range('E41').Value = '  -5.49%  '
$ws.Range('E42').Value = '  -9.53%  '
$ws.Range('E43').Value = '  -8.35%  '
Set-TextValue 'D44' '0.0562'
$ws.Range('E44').Value = '  -4.19%  '
$ws.Range('E45').Value = '  +0.16%  '
Set-TextValue 'D46' '19.66'
$ws.Range('E46').Value = '  -10.06%  '
$ws.Range('E47').Value = '  -0.63%  '
Set-TextValue 'D48' '0.0228'
$ws.Range('E48').Value = '  -5.99%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D49' '0.0890'
$ws.Range('E49').Value = '  -5.91%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D50' '4.55'
$ws.Range('E50').Value = '  -12.87%  '
Set-TextValue 'D51' '17.74'
$ws.Range('E51').Value = '  -7.90%  '
